$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Chad"
$ws.Range("C4").Value = "O"
$ws.Range("D4").Value = "demo@email.com"
$ws.Range("E4").Value = 107.6
$ws.Range("F4").Value = 19
$ws.Range("G4").Value = "Chest Pain,Loss of Movement,Fever,Dry Cough,Aches,Sore Throat"
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = $false
$ws.Range("M4").Value = "Very High Risk"
